# Generate Report for Handoff
#
# Re-running the handoff-report generator re-stamped the "Latest Handoff
# Datetime" column (D) for every row that was still awaiting a successful
# handback (status "Handback transform failed" / "Ready for handoff") on
# both the zh-cn and de-de localization-status sheets. Each locale's
# scattered 11:42:xx timestamps collapse to one fresh, later timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$newZhCnTime = "2016-03-10 11:42:56"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $zhcn.Cells.Item($r, 4).Value = $newZhCnTime
}

$dede = $wb.Worksheets.Item("de-de")
$newDeDeTime = "2016-03-10 11:43:00"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $dede.Cells.Item($r, 4).Value = $newDeDeTime
}
